$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.113642
$ws.Range("H2").Value = 0.340926
$ws.Range("I2").Value = 0.7883210666148713
$ws.Range("J2").Value = 0.7883210666148713
$ws.Range("M2").Value = 38.55267666666666
$ws.Range("N2").Value = 115.65803
$ws.Range("O2").Value = 0.5758151725879548
$ws.Range("P2").Value = 0.5758151725879548
$ws.Range("Q2").Value = 4.381203281753333
$ws.Range("R2").Value = 39.43082953578
$ws.Range("S2").Value = 0.4539272310275627
$ws.Range("T2").Value = 0.4539272310275627

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.113642
$ws.Range("H3").Value = 0.340926
$ws.Range("I3").Value = 0.7883210666148713
$ws.Range("J3").Value = 0.7883210666148713
$ws.Range("O3").Value = 0.08021535714867321
$ws.Range("P3").Value = 0.08021535714867323
$ws.Range("Q3").Value = 0.6103343619920001
$ws.Range("R3").Value = 5.493009257928001
$ws.Range("S3").Value = 0.0632354559063349
$ws.Range("T3").Value = 0.06323545590633492

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.113642
$ws.Range("H4").Value = 0.340926
$ws.Range("I4").Value = 0.7883210666148713
$ws.Range("J4").Value = 0.7883210666148713
$ws.Range("M4").Value = 23.02986166666667
$ws.Range("N4").Value = 69.089585
$ws.Range("O4").Value = 0.3439694702633719
$ws.Range("P4").Value = 0.3439694702633719
$ws.Range("Q4").Value = 2.617159539523333
$ws.Range("R4").Value = 23.55443585571
$ws.Range("S4").Value = 0.2711583796809736
$ws.Range("T4").Value = 0.2711583796809736

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.030515
$ws.Range("H5").Value = 0.091545
$ws.Range("I5").Value = 0.2116789333851287
$ws.Range("J5").Value = 0.2116789333851287
$ws.Range("M5").Value = 38.55267666666666
$ws.Range("N5").Value = 115.65803
$ws.Range("O5").Value = 0.5758151725879548
$ws.Range("P5").Value = 0.5758151725879548
$ws.Range("Q5").Value = 1.176434928483333
$ws.Range("R5").Value = 10.58791435635
$ws.Range("S5").Value = 0.1218879415603921
$ws.Range("T5").Value = 0.1218879415603921

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.030515
$ws.Range("H6").Value = 0.091545
$ws.Range("I6").Value = 0.2116789333851287
$ws.Range("J6").Value = 0.2116789333851287
$ws.Range("O6").Value = 0.08021535714867321
$ws.Range("P6").Value = 0.08021535714867323
$ws.Range("Q6").Value = 0.16388617814
$ws.Range("R6").Value = 1.47497560326
$ws.Range("S6").Value = 0.0169799012423383
$ws.Range("T6").Value = 0.01697990124233831

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.030515
$ws.Range("H7").Value = 0.091545
$ws.Range("I7").Value = 0.2116789333851287
$ws.Range("J7").Value = 0.2116789333851287
$ws.Range("M7").Value = 23.02986166666667
$ws.Range("N7").Value = 69.089585
$ws.Range("O7").Value = 0.3439694702633719
$ws.Range("P7").Value = 0.3439694702633719
$ws.Range("Q7").Value = 0.7027562287583333
$ws.Range("R7").Value = 6.324806058825
$ws.Range("S7").Value = 0.0728110905823983
$ws.Range("T7").Value = 0.07281109058239832
